$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 ("Save") should carry the same style as the other
# header cells (e.g. G1): bold, bordered, centered. Copy the formatting
# from G1 before writing the new text so the style index lines up.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
